$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = 1520.09097147339
$ws.Range("A22").Value = 68523.995
$ws.Range("A23").Value = 66852.495
$ws.Range("A24").Value = 33754.2325
